$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B5").Value = "Yahoo UK"
$ws.Range("B5").Select()
